# Update cryptocurrency price/volume figures per latest data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'24.509.77"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.25%  '
$ws.Range("D3").Value = "'1.694.28"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.15%  '
$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.36%  '
$ws.Range("D5").Value = "'313.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.56%  '
$ws.Range("D6").Value = "'1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.50%  '
$ws.Range("E7").Value = '  +0.99%  '
$ws.Range("E8").Value = '  +2.37%  '
$ws.Range("E9").Value = '  -0.39%  '
$ws.Range("D10").Value = "'1.523"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +9.59%  '
$ws.Range("D11").Value = "'53.57"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +10.83%  '
$ws.Range("D12").Value = "'0.08780"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.03%  '
$ws.Range("D13").Value = "'7.306"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +12.67%  '
$ws.Range("D14").Value = "'23.16"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.24%  '
$ws.Range("E15").Value = '  +2.89%  '
$ws.Range("D16").Value = "'7.528"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +5.79%  '
$ws.Range("D17").Value = "'1.691.52"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.44%  '
$ws.Range("D18").Value = "'100.33"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.78%  '
$ws.Range("E19").Value = '  +4.55%  '
$ws.Range("D20").Value = "'19.44"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.75%  '
$ws.Range("D21").Value = "'6.694"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.36%  '
$ws.Range("D22").Value = "'1.004"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.13%  '
$ws.Range("E23").Value = '  +3.08%  '
$ws.Range("D24").Value = "'24.492.90"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.24%  '
$ws.Range("D25").Value = "'3.008"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +9.03%  '
$ws.Range("E26").Value = '  +0.32%  '
$ws.Range("D27").Value = "'22.38"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.74%  '
$ws.Range("D28").Value = "'159.36"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.63%  '
$ws.Range("D29").Value = "'5.168"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.51%  '
$ws.Range("D30").Value = "'133.80"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.18%  '
$ws.Range("D31").Value = "'7.508"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +30.91%  '
$ws.Range("D32").Value = "'1.878.50"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.35%  '
$ws.Range("D33").Value = "'1.087"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.45%  '
$ws.Range("D34").Value = "'0.08649"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.12%  '
$ws.Range("D35").Value = "'7.340"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +19.81%  '
$ws.Range("D36").Value = "'1.971"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +8.47%  '
$ws.Range("D37").Value = "'11.04"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.96%  '
$ws.Range("D38").Value = "'0.2714"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.39%  '
$ws.Range("D39").Value = "'14.74"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.46%  '
$ws.Range("D40").Value = "'0.02761"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +10.27%  '
$ws.Range("D41").Value = "'0.08988"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.09%  '
$ws.Range("D42").Value = "'1.476"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.40%  '
$ws.Range("D43").Value = "'0.7645"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.58%  '
$ws.Range("D44").Value = "'0.7155"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.20%  '
$ws.Range("D45").Value = "'15.52"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.23%  '
$ws.Range("D46").Value = "'2.448"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.47%  '
$ws.Range("D47").Value = "'4.165"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.71%  '
$ws.Range("D48").Value = "'1.002"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.24%  '
$ws.Range("D49").Value = "'140.10"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.10%  '
$ws.Range("D50").Value = "'1.299"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +15.02%  '
$ws.Range("D51").Value = "'0.00000000378"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.40%  '
